# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" column (D) for the rows that were
# included in this handoff run, on both the zh-cn and de-de status sheets.

$wb = $excel.ActiveWorkbook

# Rows (in each localized status sheet) whose "Latest Handoff Datetime"
# (column D) is refreshed to the new handoff run's timestamp.
$handoffRows = 7,10,11,12,13,14,15,16

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$newHandoffTime_zhcn = "2016-03-10 06:37:32"
foreach ($r in $handoffRows) {
    $ws_zhcn.Cells.Item($r, 4).Value = $newHandoffTime_zhcn
}

$ws_dede = $wb.Worksheets.Item("de-de")
$newHandoffTime_dede = "2016-03-10 06:37:39"
foreach ($r in $handoffRows) {
    $ws_dede.Cells.Item($r, 4).Value = $newHandoffTime_dede
}
